# khl_stats_1369_ext.xlsx update (2025-11-01 build)
# - Adds 4 new matches (rows 414-417) to Matches_SOG
# - Refreshes as_of_utc + derived shot stats on Shots_HA / Shots_Summary
# - Bumps Meta_ext as_of_utc + build_version

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Matches_SOG: append the 4 newly played matches
# ---------------------------------------------------------------------
$wsMatches = $wb.Worksheets.Item("Matches_SOG")

$newMatches = @(
    @{ Row = 414; Uid = "897711"; Date = "2025-10-31T17:00:00"; Home = "Автомобилист"; Away = "Сибирь";        Sog = 25; SogA = 20 },
    @{ Row = 415; Uid = "897712"; Date = "2025-10-31T17:00:00"; Home = "Трактор";       Away = "Барыс";         Sog = 33; SogA = 28 },
    @{ Row = 416; Uid = "897713"; Date = "2025-10-31T17:00:00"; Home = "Салават Юлаев"; Away = "Адмирал";       Sog = 18; SogA = 39 },
    @{ Row = 417; Uid = "897664"; Date = "2025-10-31T19:30:00"; Home = "Локомотив";     Away = "СКА";           Sog = 33; SogA = 34 }
)

foreach ($m in $newMatches) {
    $r = $m.Row
    # Column A ("uid") holds numeric-looking text in every existing row, so
    # force text storage (then drop back to the Normal style so no new
    # cell format is introduced) before writing the value.
    $wsMatches.Range("A$r").NumberFormat = "@"
    $wsMatches.Range("A$r").Value = $m.Uid
    $wsMatches.Range("A$r").Style = "Normal"

    $wsMatches.Range("B$r").Value = $m.Date
    $wsMatches.Range("C$r").Value = $m.Home
    $wsMatches.Range("D$r").Value = $m.Away
    $wsMatches.Range("E$r").Value = $m.Sog
    $wsMatches.Range("F$r").Value = $m.SogA
    $wsMatches.Range("G$r").Value = "khl_text"
}

# ---------------------------------------------------------------------
# Shots_HA: bump as_of_utc for every team row, refresh the rows whose
# home/away shots-on-goal totals moved because of the new matches
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Shots_HA")

for ($r = 2; $r -le 23; $r++) {
    $ws2.Range("D$r").Value = "2025-10-31T19:30:00Z"
}

$ws2.Range("E3").Value = 17
$ws2.Range("G3").Value = 497
$ws2.Range("H3").Value = 534
$ws2.Range("I3").Value = 29.2
$ws2.Range("J3").Value = 31.4

$ws2.Range("F4").Value = 19
$ws2.Range("K4").Value = 608
$ws2.Range("L4").Value = 540
$ws2.Range("M4").Value = 32
$ws2.Range("N4").Value = 28.4

$ws2.Range("F7").Value = 15
$ws2.Range("K7").Value = 415
$ws2.Range("L7").Value = 515
$ws2.Range("M7").Value = 27.7
$ws2.Range("N7").Value = 34.3

$ws2.Range("E12").Value = 16
$ws2.Range("G12").Value = 501
$ws2.Range("H12").Value = 437
$ws2.Range("I12").Value = 31.3
$ws2.Range("J12").Value = 27.3

$ws2.Range("F15").Value = 15
$ws2.Range("K15").Value = 468
$ws2.Range("L15").Value = 487
$ws2.Range("M15").Value = 31.2
$ws2.Range("N15").Value = 32.5

$ws2.Range("E16").Value = 14
$ws2.Range("G16").Value = 387
$ws2.Range("H16").Value = 381
$ws2.Range("I16").Value = 27.6
$ws2.Range("J16").Value = 27.2

$ws2.Range("F18").Value = 20
$ws2.Range("K18").Value = 555
$ws2.Range("L18").Value = 611
$ws2.Range("M18").Value = 27.8
$ws2.Range("N18").Value = 30.6

$ws2.Range("E21").Value = 16
$ws2.Range("G21").Value = 519
$ws2.Range("H21").Value = 489
$ws2.Range("J21").Value = 30.6

# ---------------------------------------------------------------------
# Shots_Summary: bump as_of_utc for every team row, refresh the rows
# whose totals moved because of the new matches
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("Shots_Summary")

for ($r = 2; $r -le 23; $r++) {
    $ws3.Range("D$r").Value = "2025-10-31T19:30:00Z"
}

$ws3.Range("E3").Value = 41
$ws3.Range("F3").Value = 1165
$ws3.Range("G3").Value = 1259
$ws3.Range("H3").Value = 28.4
$ws3.Range("I3").Value = 30.7

$ws3.Range("E4").Value = 33
$ws3.Range("F4").Value = 1152
$ws3.Range("G4").Value = 932
$ws3.Range("H4").Value = 34.9
$ws3.Range("I4").Value = 28.2

$ws3.Range("E7").Value = 41
$ws3.Range("F7").Value = 1227
$ws3.Range("G7").Value = 1341
$ws3.Range("H7").Value = 29.9

$ws3.Range("E12").Value = 40
$ws3.Range("F12").Value = 1241
$ws3.Range("G12").Value = 1029
$ws3.Range("I12").Value = 25.7

$ws3.Range("E15").Value = 38
$ws3.Range("F15").Value = 1236
$ws3.Range("G15").Value = 1260

$ws3.Range("E16").Value = 38
$ws3.Range("F16").Value = 1045
$ws3.Range("G16").Value = 1079
$ws3.Range("H16").Value = 27.5
$ws3.Range("I16").Value = 28.4

$ws3.Range("E18").Value = 38
$ws3.Range("F18").Value = 1041
$ws3.Range("G18").Value = 1305
$ws3.Range("H18").Value = 27.4
$ws3.Range("I18").Value = 34.3

$ws3.Range("E21").Value = 40
$ws3.Range("F21").Value = 1345
$ws3.Range("G21").Value = 1273
$ws3.Range("I21").Value = 31.8

# ---------------------------------------------------------------------
# Meta_ext: bump as_of_utc + build_version
# ---------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("Meta_ext")
$ws4.Range("B2").Value = "2025-10-31T19:30:00Z"
$ws4.Range("D2").Value = 29
